$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.258.11'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '1.891.57'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'322.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.25%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("D8").Value = "'0.4040"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.67%  '
$ws.Range("D9").Value = "'47.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D10").Value = "'0.08017"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("D11").Value = "'0.9958"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.48%  '
$ws.Range("D12").Value = "'23.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.45%  '
$ws.Range("D13").Value = '1.887.22'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").Value = "'5.925"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = "'7.029"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.68%  '
$ws.Range("D16").Value = "'89.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.06631"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").Value = "'17.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").Value = '29.271.52'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = "'5.493"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").Value = "'11.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.68%  '
$ws.Range("D25").Value = "'2.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("D26").Value = '2.119.47'
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").Value = "'154.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").Value = "'19.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.11%  '
$ws.Range("D29").Value = "'5.921"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.69%  '
$ws.Range("D30").Value = "'2.077"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.56%  '
$ws.Range("D31").Value = "'117.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = "'1.031"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.97%  '
$ws.Range("D33").Value = "'0.09424"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'3.529"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("D36").Value = "'5.344"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("D37").Value = "'0.02247"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("D38").Value = "'0.06041"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("D39").Value = "'1.169"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.08%  '
$ws.Range("D40").Value = "'7.926"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.11%  '
$ws.Range("D41").Value = "'0.5825"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").Value = "'0.1832"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = "'10.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("D44").Value = "'1.289"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.40%  '
$ws.Range("D45").Value = "'0.07702"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.40%  '
$ws.Range("D46").Value = "'2.360"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("D47").Value = "'12.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = "'0.5473"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").Value = "'1.907"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'113.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = "'0.2940"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.17%  '
